$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, matching style of existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-51
$data = @{
    2 = @(6, 7)
    3 = @(7, 7)
    4 = @(7, 7)
    5 = @(1, 1)
    6 = @(7, 7)
    7 = @(5, 5)
    8 = @(6, 7)
    9 = @(7, 7)
    10 = @(7, 7)
    11 = @(3, 3)
    12 = @(1, 2)
    13 = @(8, 8)
    14 = @(7, 7)
    15 = @(7, 7)
    16 = @(1, 2)
    17 = @(10, 10)
    18 = @(4, 5)
    19 = @(1, 1)
    20 = @(1, 2)
    21 = @(7, 8)
    22 = @(1, 1)
    23 = @(7, 7)
    24 = @(8, 8)
    25 = @(7, 9)
    26 = @(1, 3)
    27 = @(1, 3)
    28 = @(8, 8)
    29 = @(5, 7)
    30 = @(7, 8)
    31 = @(8, 8)
    32 = @(9, 9)
    33 = @(1, 2)
    34 = @(8, 8)
    35 = @(1, 1)
    36 = @(9, 9)
    37 = @(1, 1)
    38 = @(1, 1)
    39 = @(1, 2)
    40 = @(2, 3)
    41 = @(8, 8)
    42 = @(4, 6)
    43 = @(7, 7)
    44 = @(2, 2)
    45 = @(7, 7)
    46 = @(5, 6)
    47 = @(5, 6)
    48 = @(8, 8)
    49 = @(6, 7)
    50 = @(6, 7)
    51 = @(2, 2)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
